$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: int/keyword/<INT_TK> -> char/keyword/<CHAR_TK>
$ws.Range("A2").Value = "char"
$ws.Range("C2").Value = "<CHAR_TK>"

# Row 3: main/keyword/<MAIN_TK> -> valid/identifier/<ID_TK, 1>
$ws.Range("A3").Value = "valid"
$ws.Range("B3").Value = "identifier"
$ws.Range("C3").Value = "<ID_TK, 1>"

# Row 7: int/keyword/<INT_TK> -> string/keyword/<STRING_TK>
$ws.Range("A7").Value = "string"
$ws.Range("C7").Value = "<STRING_TK>"

# Row 8: 1stPlace/identifier/<ID_TK, 1> -> itisvalid/identifier/<ID_TK, 2>
$ws.Range("A8").Value = "itisvalid"
$ws.Range("C8").Value = "<ID_TK, 2>"

# Row 10: 100/integer/<INT_CONST> -> 4/integer/<INT_CONST> (keep as text, not number)
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "4"
$ws.Range("A10").Style = "Normal"

# Row 12: return/keyword/<RETURN_TK> -> printf/keyword/<PRINTF_TK>
$ws.Range("A12").Value = "printf"
$ws.Range("C12").Value = "<PRINTF_TK>"

# Row 13: 0/integer/<INT_CONST> -> (/symbol/<PHARANTESES1_TK>
$ws.Range("A13").Value = "("
$ws.Range("B13").Value = "symbol"
$ws.Range("C13").Value = "<PHARANTESES1_TK>"

# Row 14: ;/symbol/<SEMICOLON_TK> -> "again it is valid"/string_constant/<STR_TK>
$ws.Range("A14").Value = '"again it is valid"'
$ws.Range("B14").Value = "string_constant"
$ws.Range("C14").Value = "<STR_TK>"

# Row 15 (new): )/symbol/<PHARANTESES2_TK>
$ws.Range("A15").Value = ")"
$ws.Range("B15").Value = "symbol"
$ws.Range("C15").Value = "<PHARANTESES2_TK>"

# Row 16 (new): printf/keyword/<PRINTF_TK>
$ws.Range("A16").Value = "printf"
$ws.Range("B16").Value = "keyword"
$ws.Range("C16").Value = "<PRINTF_TK>"

# Row 17 (new): (/symbol/<PHARANTESES1_TK>
$ws.Range("A17").Value = "("
$ws.Range("B17").Value = "symbol"
$ws.Range("C17").Value = "<PHARANTESES1_TK>"

# Row 18 (new): "/*i want to find bug in the workkk*/"/string_constant/<STR_TK>
$ws.Range("A18").Value = '"/*i want to find bug in the workkk*/"'
$ws.Range("B18").Value = "string_constant"
$ws.Range("C18").Value = "<STR_TK>"

# Row 19 (new): )/symbol/<PHARANTESES2_TK>
$ws.Range("A19").Value = ")"
$ws.Range("B19").Value = "symbol"
$ws.Range("C19").Value = "<PHARANTESES2_TK>"

# Row 20 (new): return/keyword/<RETURN_TK>
$ws.Range("A20").Value = "return"
$ws.Range("B20").Value = "keyword"
$ws.Range("C20").Value = "<RETURN_TK>"

# Row 21 (new): "lexical don't give error"/string_constant/<STR_TK>
$ws.Range("A21").Value = '"lexical don''t give error"'
$ws.Range("B21").Value = "string_constant"
$ws.Range("C21").Value = "<STR_TK>"

# Row 22: } / symbol / <BRACKET2_TK> (moved from old row 15)
$ws.Range("A22").Value = "}"
$ws.Range("B22").Value = "symbol"
$ws.Range("C22").Value = "<BRACKET2_TK>"
